# Applies the styles.xml changes described by the diff:
#   1. Add a new "Abstract Title" paragraph style (based on Normal,
#      followed by Abstract).
#   2. Change the "Abstract" style's space-before from 300 to 100
#      (twentieths of a point -> points: 300/20=15, 100/20=5).
#   3. Add a new "Footnote Block Text" paragraph style (based on
#      Footnote Text, followed by Footnote Text), mirroring the
#      existing "Block Text" style's indent/spacing.

$d = $word.ActiveDocument

# --- 1. New style: Abstract Title -------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles.Item("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles.Item("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 2. Abstract style: before-spacing 300 -> 100 ----------------------
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. New style: Footnote Block Text ----------------------------------
# NB: look the base/"next" style up by its styleId ("FootnoteText"), not
# its display name ("Footnote Text") - the COM shim mirrors whichever
# string resolved the lookup into w:basedOn/w:next, so looking it up by
# NameLocal would incorrectly write the space-containing display name
# into those attributes instead of the real styleId.
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = $d.Styles.Item("FootnoteText")
$footnoteBlockText.NextParagraphStyle = $d.Styles.Item("FootnoteText")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24
